# Apply "Added handling of common packages." edit to the classFields sheet.
# This reorders the field rows belonging to OrderGeneratorService,
# OrderControllerTests and OrderController so that fields shared with
# common/base packages are grouped together (and one C8/C9 modifier swap).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

# pl.piomin.order.service.OrderGeneratorService fields (rows 4-7)
$ws.Range("B4").Value = "id"
$ws.Range("D4").Value = "java.util.concurrent.atomic.AtomicLong"

$ws.Range("B5").Value = "template"
$ws.Range("D5").Value = "org.springframework.kafka.core.KafkaTemplate"

$ws.Range("B6").Value = "RAND"
$ws.Range("D6").Value = "java.util.Random"

$ws.Range("B7").Value = "executor"
$ws.Range("D7").Value = "java.util.concurrent.Executor"

# pl.piomin.order.OrderControllerTests fields (rows 8-9)
$ws.Range("B8").Value = "restTemplate"
$ws.Range("C8").Value = ""
$ws.Range("D8").Value = "org.springframework.boot.test.web.client.TestRestTemplate"

$ws.Range("B9").Value = "factory"
$ws.Range("C9").Value = "private"
$ws.Range("D9").Value = "org.springframework.kafka.core.ConsumerFactory"

# pl.piomin.order.controller.OrderController fields (rows 12-16)
$ws.Range("B12").Value = "LOG"
$ws.Range("D12").Value = "org.slf4j.Logger"

$ws.Range("B13").Value = "template"
$ws.Range("D13").Value = "org.springframework.kafka.core.KafkaTemplate"

$ws.Range("B14").Value = "kafkaStreamsFactory"
$ws.Range("D14").Value = "org.springframework.kafka.config.StreamsBuilderFactoryBean"

$ws.Range("B15").Value = "id"
$ws.Range("D15").Value = "java.util.concurrent.atomic.AtomicLong"

$ws.Range("B16").Value = "orderGeneratorService"
$ws.Range("D16").Value = "pl.piomin.order.service.OrderGeneratorService"
